$d = $word.ActiveDocument

$pairs = @(
    @("339÷9=37, 6", "115÷7=16, 3"),
    @("766÷3=255, 1", "931÷4=232, 3"),
    @("355÷5=71, 0", "613÷5=122, 3"),
    @("713÷9=79, 2", "580÷6=96, 4"),
    @("399÷7=57, 0", "230÷6=38, 2"),
    @("927÷2=463, 1", "628÷9=69, 7"),
    @("354÷8=44, 2", "886÷5=177, 1"),
    @("301÷7=43, 0", "670÷6=111, 4"),
    @("105÷3=35, 0", "373÷7=53, 2"),
    @("813÷7=116, 1", "260÷3=86, 2"),
    @("157÷7=22, 3", "633÷8=79, 1"),
    @("606÷2=303, 0", "855÷6=142, 3"),
    @("134÷7=19, 1", "118÷8=14, 6"),
    @("141÷6=23, 3", "562÷7=80, 2"),
    @("453÷7=64, 5", "211÷2=105, 1"),
    @("623÷6=103, 5", "960÷9=106, 6"),
    @("759÷9=84, 3", "445÷5=89, 0"),
    @("661÷6=110, 1", "393÷4=98, 1"),
    @("887÷9=98, 5", "166÷5=33, 1"),
    @("107÷4=26, 3", "898÷9=99, 7"),
    @("289÷7=41, 2", "816÷8=102, 0"),
    @("881÷3=293, 2", "196÷7=28, 0"),
    @("295÷8=36, 7", "898÷7=128, 2"),
    @("620÷9=68, 8", "985÷7=140, 5"),
    @("113÷2=56, 1", "888÷4=222, 0")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
